$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry")

# Insert a new column before AN, shifting thickness_max_chord_ratio (and
# everything after it) one column to the right.
$ws.Columns("AN:AN").Insert()

# Populate the newly inserted column with the "solidity" data.
$ws.Range("AN1").Value = "solidity"
$ws.Range("AN2").Value = "[1.42997704 1.70997375]"
